# Update column F (dSF) values in several rows on Sheet1.
# These edits correspond to a data repull / recalculation of the
# "mean" derived dSF column that previously mirrored dS0 (column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -8
    3  = -6
    8  = -7
    11 = -4
    13 = -5
    15 = -3
    18 = -3
    22 = -2
    24 = 6
    29 = -2
    30 = -1
    31 = -3
    33 = -5
    37 = 2
    38 = 13
    41 = -4
    42 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
